$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "fixed horizontal centering on registers" -> enable "Center on page: Horizontally"
# in Page Setup (Print Options), matching the new <printOptions horizontalCentered="1"/>
$ws.PageSetup.CenterHorizontally = $true

# Move the active selection from A5:M5 to the single cell G12,
# matching the updated <selection activeCell="G12" sqref="G12"/>.
$null = $ws.Range("G12").Select()
